# chore: update Sheets via scheduled runner
# Refresh cached market-board figures (currentAveragePrice / LevePrice* /
# LeveProfit* columns H-N) for the affected leve rows on each class sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 111.21739
$ws.Range("J55").Value = 106.63636
$ws.Range("L55").Value = 106.63636
$ws.Range("N55").Value = -534.63636
$ws.Range("H70").Value = 4889.143
$ws.Range("I70").Value = 4000
$ws.Range("J70").Value = 4957.5386
$ws.Range("K70").Value = 12000
$ws.Range("L70").Value = 14872.6158
$ws.Range("M70").Value = -11730
$ws.Range("N70").Value = -15412.6158
$ws.Range("H73").Value = 4889.143
$ws.Range("I73").Value = 4000
$ws.Range("J73").Value = 4957.5386
$ws.Range("K73").Value = 12000
$ws.Range("L73").Value = 14872.6158
$ws.Range("M73").Value = -11064
$ws.Range("N73").Value = -16744.6158
$ws.Range("H86").Value = 6626.4443
$ws.Range("I86").Value = 6533.6924
$ws.Range("K86").Value = 6533.6924
$ws.Range("M86").Value = -5410.6924
$ws.Range("H89").Value = 6626.4443
$ws.Range("I89").Value = 6533.6924
$ws.Range("K89").Value = 32668.462
$ws.Range("M89").Value = -27052.462
$ws.Range("H130").Value = 69950
$ws.Range("J130").Value = 69950
$ws.Range("L130").Value = 69950
$ws.Range("N130").Value = -79990
$ws.Range("H132").Value = 30013
$ws.Range("I132").Value = 30013
$ws.Range("K132").Value = 90039
$ws.Range("M132").Value = -87509
$ws.Range("H141").Value = 10638.8
$ws.Range("I141").Value = 3333
$ws.Range("J141").Value = 13295.454
$ws.Range("K141").Value = 9999
$ws.Range("L141").Value = 39886.362
$ws.Range("M141").Value = -4819
$ws.Range("N141").Value = -50246.362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 50000
$ws.Range("J24").Value = 50000
$ws.Range("L24").Value = 50000
$ws.Range("N24").Value = -50748
$ws.Range("H45").Value = 1343.7693
$ws.Range("I45").Value = 977.5
$ws.Range("J45").Value = 1929.8
$ws.Range("K45").Value = 977.5
$ws.Range("L45").Value = 1929.8
$ws.Range("M45").Value = -600.5
$ws.Range("N45").Value = -2683.8
$ws.Range("H61").Value = 2779761.8
$ws.Range("I61").Value = 2040.5428
$ws.Range("K61").Value = 2040.5428
$ws.Range("M61").Value = -1828.5428
$ws.Range("H100").Value = 50000
$ws.Range("J100").Value = 50000
$ws.Range("L100").Value = 50000
$ws.Range("N100").Value = -52164
$ws.Range("H132").Value = 4156.1
$ws.Range("I132").Value = 2306.3333
$ws.Range("K132").Value = 6918.999899999999
$ws.Range("M132").Value = -4388.999899999999
$ws.Range("H134").Value = 75000.5
$ws.Range("J134").Value = 75000.5
$ws.Range("L134").Value = 75000.5
$ws.Range("N134").Value = -85140.5
$ws.Range("H136").Value = 2779761.8
$ws.Range("I136").Value = 2040.5428
$ws.Range("K136").Value = 6121.6284
$ws.Range("M136").Value = -3571.6284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3789409
$ws.Range("I134").Value = 1383.15
$ws.Range("K134").Value = 4149.450000000001
$ws.Range("M134").Value = -1614.450000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 148.52942
$ws.Range("I7").Value = 127.454544
$ws.Range("K7").Value = 127.454544
$ws.Range("M7").Value = -14.454544
$ws.Range("H28").Value = 83994.8
$ws.Range("J28").Value = 29993.5
$ws.Range("L28").Value = 29993.5
$ws.Range("N28").Value = -30483.5
$ws.Range("H86").Value = 42473.832
$ws.Range("I86").Value = 395894.5
$ws.Range("K86").Value = 395894.5
$ws.Range("M86").Value = -394771.5
$ws.Range("H89").Value = 42473.832
$ws.Range("I89").Value = 395894.5
$ws.Range("K89").Value = 1979472.5
$ws.Range("M89").Value = -1973856.5
$ws.Range("H132").Value = 2359.3809
$ws.Range("I132").Value = 2057.389
$ws.Range("K132").Value = 6172.167
$ws.Range("M132").Value = -3642.167
$ws.Range("H134").Value = 4623.8213
$ws.Range("J134").Value = 5754.5454
$ws.Range("L134").Value = 17263.6362
$ws.Range("N134").Value = -22333.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 139.38461
$ws.Range("J2").Value = 227.14285
$ws.Range("L2").Value = 1362.8571
$ws.Range("N2").Value = -1588.8571
$ws.Range("H14").Value = 400
$ws.Range("I14").Value = 400
$ws.Range("K14").Value = 1200
$ws.Range("M14").Value = -1027
$ws.Range("H34").Value = 935.55554
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 76.25
$ws.Range("I2").Value = 53.5
$ws.Range("K2").Value = 53.5
$ws.Range("M2").Value = 59.5
$ws.Range("H70").Value = 8312.145500000001
$ws.Range("I70").Value = 8555.275
$ws.Range("J70").Value = 7941.0527
$ws.Range("K70").Value = 8555.275
$ws.Range("L70").Value = 7941.0527
$ws.Range("M70").Value = -8285.275
$ws.Range("N70").Value = -8481.0527
$ws.Range("H73").Value = 8312.145500000001
$ws.Range("I73").Value = 8555.275
$ws.Range("J73").Value = 7941.0527
$ws.Range("K73").Value = 8555.275
$ws.Range("L73").Value = 7941.0527
$ws.Range("M73").Value = -7619.275
$ws.Range("N73").Value = -9813.0527
$ws.Range("H92").Value = 10000
$ws.Range("J92").Value = 10000
$ws.Range("L92").Value = 10000
$ws.Range("N92").Value = -13744
$ws.Range("H102").Value = 1551.2307
$ws.Range("I102").Value = 1528.6
$ws.Range("K102").Value = 1528.6
$ws.Range("M102").Value = 93.40000000000009
$ws.Range("H126").Value = 1990
$ws.Range("I126").Value = 1050
$ws.Range("J126").Value = 2225
$ws.Range("K126").Value = 3150
$ws.Range("L126").Value = 6675
$ws.Range("M126").Value = -680
$ws.Range("N126").Value = -11615
$ws.Range("H132").Value = 13552.518
$ws.Range("I132").Value = 7360.96
$ws.Range("K132").Value = 22082.88
$ws.Range("M132").Value = -19552.88

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2721.1614
$ws.Range("I22").Value = 2305.2307
$ws.Range("K22").Value = 2305.2307
$ws.Range("M22").Value = -2010.2307
$ws.Range("H23").Value = 20000
$ws.Range("I23").Value = 20000
$ws.Range("K23").Value = 20000
$ws.Range("M23").Value = -19770
$ws.Range("H27").Value = 2721.1614
$ws.Range("I27").Value = 2305.2307
$ws.Range("K27").Value = 2305.2307
$ws.Range("M27").Value = -2198.2307
$ws.Range("H82").Value = 2768.923
$ws.Range("I82").Value = 3199.5
$ws.Range("J82").Value = 2399.8572
$ws.Range("K82").Value = 3199.5
$ws.Range("L82").Value = 2399.8572
$ws.Range("M82").Value = -2838.5
$ws.Range("N82").Value = -3121.8572
$ws.Range("H85").Value = 2768.923
$ws.Range("I85").Value = 3199.5
$ws.Range("J85").Value = 2399.8572
$ws.Range("K85").Value = 3199.5
$ws.Range("L85").Value = 2399.8572
$ws.Range("M85").Value = -1951.5
$ws.Range("N85").Value = -4895.8572
$ws.Range("H122").Value = 3400.0557
$ws.Range("J122").Value = 4870.5
$ws.Range("L122").Value = 14611.5
$ws.Range("N122").Value = -19511.5
$ws.Range("H132").Value = 1518386.8
$ws.Range("I132").Value = 4169471
$ws.Range("J132").Value = 3481.5715
$ws.Range("K132").Value = 12508413
$ws.Range("L132").Value = 10444.7145
$ws.Range("M132").Value = -12505883
$ws.Range("N132").Value = -15504.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I62").Value = 33263.332
$ws.Range("J62").Value = 9750
$ws.Range("K62").Value = 33263.332
$ws.Range("L62").Value = 9750
$ws.Range("M62").Value = -32639.332
$ws.Range("N62").Value = -10998
$ws.Range("I65").Value = 33263.332
$ws.Range("J65").Value = 9750
$ws.Range("K65").Value = 166316.66
$ws.Range("L65").Value = 48750
$ws.Range("M65").Value = -163196.66
$ws.Range("N65").Value = -54990
$ws.Range("H122").Value = 49111.082
$ws.Range("J122").Value = 226752.2
$ws.Range("L122").Value = 680256.6000000001
$ws.Range("N122").Value = -685156.6000000001
$ws.Range("H131").Value = 71500
$ws.Range("J131").Value = 71500
$ws.Range("L131").Value = 71500
$ws.Range("N131").Value = -81580
$ws.Range("H136").Value = 5631798
$ws.Range("I136").Value = 2559068.5
$ws.Range("J136").Value = 31750000
$ws.Range("K136").Value = 7677205.5
$ws.Range("L136").Value = 95250000
$ws.Range("M136").Value = -7674655.5
$ws.Range("N136").Value = -95255100
